$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells D and E are stored as text in the source data (inline strings),
# so force text number-format before assigning to avoid Excel auto-converting
# numeric-looking strings (e.g. "1.00", "38.50") into numbers and losing
# trailing zeros / formatting.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.702.08'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.44%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.553.95'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.36%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.01'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.46'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.98%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.554.88'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.41%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.58%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.356'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.92'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.021.42'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '70.563.09'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.61%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -5.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.48'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.557.59'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.90'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.81%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -4.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '356.37'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -4.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.94'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.73%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.01'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.38'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.07'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.70%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.667.18'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.95'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.28'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '473.18'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.78'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '159.03'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.08'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.88%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.76%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.21%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.83%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -6.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '38.50'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.95%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '145.41'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.540'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.36%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.30%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.95%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0741'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.16%  '
